# Update the "Region" column (column O) values in Sheet1 of the beach/station
# data workbook. Rows whose Region was "Other" are being reclassified into
# more specific regions: "Island", "East Coast", or "West Coast".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows that move from "Other" -> "Island"
# (Hawaii, Marshall Islands, Wake Island, Bermuda, Puerto Rico stations)
$islandRows = @(2, 3, 4, 5, 6, 7, 11, 12, 13, 14, 145, 146)
foreach ($r in $islandRows) {
    $ws.Cells.Item($r, 15).Value = "Island"
}

# Rows that move from "Other" -> "East Coast"
# (Philadelphia, PA and Washington, DC stations)
$eastCoastRows = @(35, 45)
foreach ($r in $eastCoastRows) {
    $ws.Cells.Item($r, 15).Value = "East Coast"
}

# Rows that move from "Other" -> "West Coast"
# (Alaska stations)
$westCoastRows = 125..142
foreach ($r in $westCoastRows) {
    $ws.Cells.Item($r, 15).Value = "West Coast"
}

$wb.Save()
